$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the "duplicate_image_filename" column (E) with "NA" for the
# practice rows (2-5) and the main trial rows (6-21).
foreach ($r in 2..21) {
    $ws.Cells.Item($r, 5).Value = "NA"
}
